$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing comment text/author on I3 before moving it to I2
$oldComment = $ws.Range("I3").Comment

# Put the original "Estimated relationship" formula into the new I2 cell
$ws.Range("I2").Formula = '="Estimated relationship = GPA(fitted) = "&TEXT(INTERCEPT(C4:C11,D4:D11),"0.000")&"+"&TEXT(SLOPE(C4:C11,D4:D11),"0.000")&"ACT"'

# Replace I3 with the new "When ACT=20" formula
$ws.Range("I3").Formula = '="When ACT=20 the GPA will be = "&TEXT(INTERCEPT(C4:C11,D4:D11)+SLOPE(C4:C11,D4:D11)*20,"0.000")'

# Move the cell comment from I3 to I2
if ($oldComment -ne $null) {
    $oldComment.Delete()
}
$ws.Range("I2").AddComment("rodri:`nThe direction of the relationship is positive, i.e. when ACT grows (one point), GPA is expected to grow by about 0.10 score points.  `n`nThe intercept shows the GPA value when the other coefficients are equal to zero. That is, if any student in the class gets a score of 0 we expect GPA = 0.568.`n`nIf the ACT score is increased by 5 points the GPA will be about 0.51 (=0.102*5) greater.")

# Update the sheet view: drop the F1 top-left freeze and move the selection
$ws.Range("I20").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
